$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new day 15 entry (row 16): title, text, person, image_url
$ws.Range("B16").Value = "Relax am Abend - nicht müssen"
$ws.Range("C16").Value = "Ich war nach der Arbeit wieder etwas gereizt. Ich habe LS gemacht und danach wollte ich noch weitermachen mit den Coachinginhalten bzw. meinem Zukunftsbild. Aber ich war einfach nicht in Stimmung dafür und habe mir dann bewusst eine Auszeit gegönnt und einfach Darts-WM geschaut - ohne schlechtes Gewissen. "
$ws.Range("D16").Value = "frei sein"
$ws.Range("E16").Value = "https://www.dartn.de/wp-content/uploads/2025/12/Hintergrund_495-1200x800.jpg"

# Update the selected cell to match the saved workbook state
$ws.Range("D16").Select()
